# Insert a new weekly price record for "Poroto verde" / Terminal La Palmera
# de La Serena right above the current row 278, shifting the existing
# rows 278-334 down by one (row 334 -> row 335), matching the new weekly
# "Fruta / hortaliza, semanal" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 278..334 down one row, leaving row 278 free for the new record.
$ws.Rows.Item(278).Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(278, 1).Value  = 8
$ws.Cells.Item(278, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(278, 3).Value  = "Coquimbo"
$ws.Cells.Item(278, 4).Value  = 44951
$ws.Cells.Item(278, 5).Value  = 4
$ws.Cells.Item(278, 6).Value  = 100112031
$ws.Cells.Item(278, 7).Value  = "Poroto verde"
$ws.Cells.Item(278, 8).Value  = "Magnum"
$ws.Cells.Item(278, 9).Value  = "Primera"
$ws.Cells.Item(278, 10).Value = 400
$ws.Cells.Item(278, 11).Value = 23500
$ws.Cells.Item(278, 12).Value = 24000
$ws.Cells.Item(278, 13).Value = 23750
$ws.Cells.Item(278, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(278, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(278, 16).Value = 950
$ws.Cells.Item(278, 17).Value = 25
$ws.Cells.Item(278, 18).Value = "Hortaliza"
